$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

$ws1.Range("F2").Value = 136
$ws1.Range("F3").Value = 368
$ws1.Range("F4").Value = 196
$ws1.Range("F6").Value = 1236
$ws1.Range("F7").Value = 446
$ws1.Range("F8").Value = 101
$ws1.Range("F9").Value = 182
$ws1.Range("F10").Value = 150
$ws1.Range("F11").Value = 170
$ws1.Range("F12").Value = 1044
$ws1.Range("F14").Value = 269
$ws1.Range("F15").Value = 184
$ws1.Range("F16").Value = 1488
$ws1.Range("F18").Value = 224
$ws1.Range("F19").Value = 345
$ws1.Range("F21").Value = 809
$ws1.Range("F22").Value = 1147
$ws1.Range("F25").Value = 2644
$ws1.Range("F26").Value = 1421
$ws1.Range("F27").Value = 61
$ws1.Range("F28").Value = 29
$ws1.Range("F29").Value = 387
$ws1.Range("F30").Value = 413
$ws1.Range("F31").Value = 1196
$ws1.Range("F32").Value = 816
$ws1.Range("F33").Value = 1321
$ws1.Range("F34").Value = 156
$ws1.Range("F36").Value = 776
$ws1.Range("F37").Value = 586
$ws1.Range("F38").Value = 664
$ws1.Range("F39").Value = 831
$ws1.Range("F40").Value = 357
$ws1.Range("F41").Value = 242

$ws2.Range("F15").Value = 621
$ws2.Range("F22").Value = 18
$ws2.Range("F23").Value = 14

$ws4.Range("F5").Value = 136
$ws4.Range("F6").Value = 368
$ws4.Range("F7").Value = 196
$ws4.Range("F10").Value = 1236
$ws4.Range("F11").Value = 446
$ws4.Range("F12").Value = 101
$ws4.Range("F13").Value = 182
$ws4.Range("F15").Value = 150
$ws4.Range("F16").Value = 170
$ws4.Range("F18").Value = 269
$ws4.Range("F20").Value = 184
$ws4.Range("F21").Value = 1488
$ws4.Range("F23").Value = 224
$ws4.Range("F24").Value = 345
$ws4.Range("F26").Value = 1147
$ws4.Range("F27").Value = 2644
$ws4.Range("F29").Value = 1421
$ws4.Range("F30").Value = 61
$ws4.Range("F32").Value = 29
$ws4.Range("F34").Value = 387
$ws4.Range("F35").Value = 413
$ws4.Range("F36").Value = 1196
$ws4.Range("F39").Value = 816
$ws4.Range("F40").Value = 1321
$ws4.Range("F41").Value = 776
$ws4.Range("F42").Value = 586
$ws4.Range("F43").Value = 664
$ws4.Range("F44").Value = 831
$ws4.Range("F45").Value = 357
$ws4.Range("F46").Value = 18
$ws4.Range("F47").Value = 14
$ws4.Range("F48").Value = 242
